# Semana 49 de 2025
# Adds a new "week 49" column (AZ) to the weekly IRA/UCI tracking sheet,
# mirroring the existing week-number header row and filling in the
# reported counts for that week. A couple of previously-sparse cells
# (AY36, X43, AA43) are also populated as part of this week's data load.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: new week number "49" in AZ1, matching the style of the
#     existing week-number headers (e.g. AY1 = "48"). Using a leading
#     apostrophe forces the value to be stored as text (consistent with
#     the other week headers, which are text, not numbers), then we copy
#     the formatting from AY1 so AZ1 ends up with identical styling.
$ws.Range("AZ1").Value = "'49"
$ws.Range("AY1").Copy() | Out-Null
$ws.Range("AZ1").PasteSpecial(-4122) | Out-Null

# --- Data rows: week 49 counts for each facility (column AZ) ---
$ws.Range("AZ2").Value = 0
$ws.Range("AZ5").Value = 0
$ws.Range("AZ6").Value = 1
$ws.Range("AZ7").Value = 0
$ws.Range("AZ8").Value = 0
$ws.Range("AZ9").Value = 0
$ws.Range("AZ10").Value = 0
$ws.Range("AZ11").Value = 0
$ws.Range("AZ12").Value = 0
$ws.Range("AZ14").Value = 0
$ws.Range("AZ16").Value = 0
$ws.Range("AZ17").Value = 0
$ws.Range("AZ23").Value = 0
$ws.Range("AZ25").Value = 0
$ws.Range("AZ26").Value = 0
$ws.Range("AZ28").Value = 0
$ws.Range("AZ29").Value = 1
$ws.Range("AZ31").Value = 0
$ws.Range("AZ35").Value = 0

# Row 36 previously had no value at all in column AY (week 48); it is
# filled in now alongside the new AZ36 (week 49) value.
$ws.Range("AY36").Value = 1
$ws.Range("AZ36").Value = 0

$ws.Range("AZ37").Value = 0
$ws.Range("AZ38").Value = 0
$ws.Range("AZ41").Value = 0
$ws.Range("AZ42").Value = 0

# Row 43 previously had gaps at X43 and AA43 (week 21 and week 24); these
# are filled in with 0 along with the new AZ43 (week 49) value.
$ws.Range("X43").Value = 0
$ws.Range("AA43").Value = 0
$ws.Range("AZ43").Value = 0

$ws.Range("AZ44").Value = 0
$ws.Range("AZ45").Value = 0
$ws.Range("AZ46").Value = 0
$ws.Range("AZ47").Value = 0
$ws.Range("AZ48").Value = 0
$ws.Range("AZ49").Value = 0
$ws.Range("AZ50").Value = 0
$ws.Range("AZ51").Value = 0
$ws.Range("AZ54").Value = 0
$ws.Range("AZ55").Value = 0
$ws.Range("AZ56").Value = 0
$ws.Range("AZ57").Value = 0
$ws.Range("AZ58").Value = 0
$ws.Range("AZ59").Value = 0
